$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '69.141.92'
Set-TextValue "E2" '  +0.77%  '

# Row 3
Set-TextValue "D3" '2.477.86'
Set-TextValue "E3" '  +0.85%  '

# Row 4
Set-TextValue "E4" '  -0.08%  '

# Row 5
Set-TextValue "E5" '  -0.60%  '

# Row 6
Set-TextValue "D6" '162.55'
Set-TextValue "E6" '  -0.44%  '

# Row 7
Set-TextValue "E7" '  -0.11%  '

# Row 8
Set-TextValue "D8" '0.507'
Set-TextValue "E8" '  +0.09%  '

# Row 9
Set-TextValue "E9" '  +0.24%  '

# Row 10
Set-TextValue "E10" '  +0.55%  '

# Row 11
Set-TextValue "E11" '  -2.74%  '

# Row 12
Set-TextValue "D12" '4.88'
Set-TextValue "E12" '  +1.14%  '

# Row 13
Set-TextValue "E13" '  -0.10%  '

# Row 14
Set-TextValue "D14" '68.984.34'
Set-TextValue "E14" '  +0.73%  '

# Row 15
Set-TextValue "E15" '  -1.50%  '

# Row 16
Set-TextValue "D16" '23.70'
Set-TextValue "E16" '  +0.16%  '

# Row 17
Set-TextValue "D17" '2.462.73'
Set-TextValue "E17" '  -0.62%  '

# Row 18
Set-TextValue "D18" '10.74'
Set-TextValue "E18" '  -2.59%  '

# Row 19
Set-TextValue "D19" '337.45'
Set-TextValue "E19" '  -1.99%  '

# Row 20
Set-TextValue "D20" '6.98'
Set-TextValue "E20" '  -3.19%  '

# Row 21
Set-TextValue "D21" '3.80'
Set-TextValue "E21" '  -0.67%  '

# Row 22
Set-TextValue "B22" 'SuiNetwork'
Set-TextValue "C22" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D22" '1.88'
Set-TextValue "E22" '  -0.46%  '

# Row 23
Set-TextValue "B23" 'Dai'
Set-TextValue "C23" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D23" '1.00'
Set-TextValue "E23" '  -0.01%  '

# Row 24
Set-TextValue "D24" '66.95'
Set-TextValue "E24" '  -1.77%  '

# Row 25
Set-TextValue "D25" '3.68'
Set-TextValue "E25" '  -2.34%  '

# Row 26
Set-TextValue "B26" 'Aptos'
Set-TextValue "C26" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D26" '8.25'
Set-TextValue "E26" '  -0.09%  '

# Row 27
Set-TextValue "B27" 'PEPE'
Set-TextValue "C27" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D27" '0.0₃0822'
Set-TextValue "E27" '  -2.66%  '

# Row 28
Set-TextValue "B28" 'InternetComputer(DFINITY)'
Set-TextValue "C28" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D28" '7.23'
Set-TextValue "E28" '  -1.31%  '

# Row 29
Set-TextValue "B29" 'FirstDigitalUSD'
Set-TextValue "C29" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  +0.01%  '

# Row 30
Set-TextValue "B30" 'Bittensor'
Set-TextValue "C30" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D30" '431.48'
Set-TextValue "E30" '  -1.11%  '

# Row 31
Set-TextValue "B31" 'Fetch.AI'
Set-TextValue "C31" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D31" '1.14'
Set-TextValue "E31" '  -3.87%  '

# Row 32
Set-TextValue "B32" 'PancakeSwap'
Set-TextValue "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '1.62'
Set-TextValue "E32" '  -3.97%  '

# Row 33
Set-TextValue "B33" 'Monero'
Set-TextValue "C33" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D33" '158.51'
Set-TextValue "E33" '  +0.99%  '

# Row 34
Set-TextValue "B34" 'WhiteBITCoin'
Set-TextValue "C34" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D34" '19.03'
Set-TextValue "E34" '  +0.08%  '

# Row 35
Set-TextValue "B35" 'Kaspa'
Set-TextValue "C35" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D35" '0.110'
Set-TextValue "E35" '  +0.06%  '

# Row 36
Set-TextValue "B36" 'USDe'
Set-TextValue "C36" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D36" '1.00'
Set-TextValue "E36" '  -0.07%  '

# Row 37
Set-TextValue "B37" 'EthereumClassic'
Set-TextValue "C37" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D37" '17.82'
Set-TextValue "E37" '  -0.56%  '

# Row 38
Set-TextValue "B38" 'PolygonEcosystemToken'
Set-TextValue "C38" 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue "D38" '0.301'
Set-TextValue "E38" '  -2.14%  '

# Row 39
Set-TextValue "B39" 'RenderToken'
Set-TextValue "C39" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D39" '4.44'
Set-TextValue "E39" '  -1.75%  '

# Row 40
Set-TextValue "B40" 'Stacks'
Set-TextValue "C40" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D40" '1.48'
Set-TextValue "E40" '  -4.16%  '

# Row 41
Set-TextValue "B41" 'ImmutableX'
Set-TextValue "C41" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D41" '1.09'
Set-TextValue "E41" '  -1.90%  '

# Row 42
Set-TextValue "B42" 'dogwifhat'
Set-TextValue "C42" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D42" '2.08'
Set-TextValue "E42" '  -1.15%  '

# Row 43
Set-TextValue "B43" 'Filecoin'
Set-TextValue "C43" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D43" '3.36'
Set-TextValue "E43" '  -0.72%  '

# Row 44
Set-TextValue "B44" 'Aave'
Set-TextValue "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '131.73'
Set-TextValue "E44" '  -2.92%  '

# Row 45
Set-TextValue "B45" 'ARBITRUM'
Set-TextValue "C45" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D45" '0.486'
Set-TextValue "E45" '  -0.99%  '

# Row 46
Set-TextValue "B46" 'Cronos'
Set-TextValue "C46" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D46" '0.0713'
Set-TextValue "E46" '  -0.85%  '

# Row 47
Set-TextValue "B47" 'Mantle'
Set-TextValue "C47" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D47" '0.565'
Set-TextValue "E47" '  -0.05%  '

# Row 48
Set-TextValue "B48" 'Stellar'
Set-TextValue "C48" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D48" '0.0914'
Set-TextValue "E48" '  -0.14%  '

# Row 49
Set-TextValue "B49" 'BitgetToken'
Set-TextValue "C49" 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
Set-TextValue "D49" '1.12'
Set-TextValue "E49" '  +0.20%  '

# Row 50
Set-TextValue "B50" 'Optimism'
Set-TextValue "C50" 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
Set-TextValue "D50" '1.40'
Set-TextValue "E50" '  -2.75%  '

# Row 51
Set-TextValue "B51" 'THORChain'
Set-TextValue "C51" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D51" '5.02'
Set-TextValue "E51" '  -7.74%  '
